$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.058715
$ws.Range("H2").Value = 0.176145
$ws.Range("M2").Value = 30.224788
$ws.Range("N2").Value = 90.674364
$ws.Range("O2").Value = 0.3247052378228209
$ws.Range("P2").Value = 0.3247052378228209
$ws.Range("Q2").Value = 1.77464842742
$ws.Range("R2").Value = 15.97183584678
$ws.Range("S2").Value = 0.3247052378228209
$ws.Range("T2").Value = 0.3247052378228209

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.058715
$ws.Range("H3").Value = 0.176145
$ws.Range("M3").Value = 20.25845733333333
$ws.Range("O3").Value = 0.2176368352473959
$ws.Range("P3").Value = 0.217636835247396
$ws.Range("Q3").Value = 1.189475322326667
$ws.Range("R3").Value = 10.70527790094
$ws.Range("S3").Value = 0.2176368352473959
$ws.Range("T3").Value = 0.217636835247396

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.058715
$ws.Range("H4").Value = 0.176145
$ws.Range("M4").Value = 12.725178
$ws.Range("N4").Value = 38.175534
$ws.Range("O4").Value = 0.1367067305427495
$ws.Range("P4").Value = 0.1367067305427495
$ws.Range("Q4").Value = 0.74715882627
$ws.Range("R4").Value = 6.724429436429999
$ws.Range("S4").Value = 0.1367067305427495
$ws.Range("T4").Value = 0.1367067305427495

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.058715
$ws.Range("H5").Value = 0.176145
$ws.Range("M5").Value = 29.87534766666667
$ws.Range("N5").Value = 89.62604300000001
$ws.Range("O5").Value = 0.3209511963870337
$ws.Range("P5").Value = 0.3209511963870337
$ws.Range("Q5").Value = 1.754131038248333
$ws.Range("R5").Value = 15.787179344235
$ws.Range("S5").Value = 0.3209511963870337
$ws.Range("T5").Value = 0.3209511963870337
